# Improved Link error reporting.
#
# The two "Syntax error in AQL expression:" diagnostic messages emitted for
# an invalid \link statement are renamed to "Invalid link statement:".
#
# We locate each message with Find (so we do not depend on fragile
# character offsets) and then assign straight onto the found Range's
# .Text property (rather than passing the replacement through
# Find.Execute's Replace argument) so that straight quote characters are
# not auto-corrected into curly/smart quotes.

$d = $word.ActiveDocument

$oldPrefix = "Syntax error in AQL expression:"
$newPrefix = "Invalid link statement:"

$needle1 = $oldPrefix + " Expression ""self. 'a reference to bookmark1'"" is invalid: missing feature access or service call"
$replacement1 = $newPrefix + " Expression ""self. 'a reference to bookmark1'"" is invalid: missing feature access or service call"

$needle2 = $oldPrefix + " Expression """" is invalid: null or empty string."
$replacement2 = $newPrefix + " Expression """" is invalid: null or empty string."

$range1 = $d.Content
if ($range1.Find.Execute($needle1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $range1.Text = $replacement1
}

$range2 = $d.Content
if ($range2.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $range2.Text = $replacement2
}
